$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column F header from "Fraction" to "Compartment"
$ws.Range("F1").Value = "Compartment"

# Replace "Endo" values in the Compartment column with "Root"
for ($r = 2; $r -le 39; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Value2 -eq "Endo") {
        $cell.Value = "Root"
    }
}

# Update the active selection to H19
$ws.Range("H19").Select()
